$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 57.666668
$ws.Range("M8").Value = -34.00000399999999
$ws.Range("K8").Value = 173.000004
$ws.Range("I8").Value = 57.666668
$ws.Range("H40").Value = 3170.5557
$ws.Range("L40").Value = 4899.9
$ws.Range("J40").Value = 4899.9
$ws.Range("N40").Value = -5249.9
$ws.Range("I94").Value = 2230.1667
$ws.Range("K94").Value = 2230.1667
$ws.Range("H94").Value = 2230.1667
$ws.Range("M94").Value = -1779.1667
$ws.Range("H103").Value = 937.0833
$ws.Range("M103").Value = -2221.1429
$ws.Range("K103").Value = 2807.1429
$ws.Range("I103").Value = 935.7143
$ws.Range("K137").Value = 0
$ws.Range("J137").Value = 3013.2856
$ws.Range("I137").Value = 0
$ws.Range("N137").Value = -14139.8568
$ws.Range("L137").Value = 9039.856800000001
$ws.Range("M137").ClearContents()
$ws.Range("H137").Value = 3013.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J32").Value = 16799.6
$ws.Range("I32").Value = 3155.6316
$ws.Range("M32").Value = -2868.6316
$ws.Range("N32").Value = -17373.6
$ws.Range("L32").Value = 16799.6
$ws.Range("H32").Value = 5998.125
$ws.Range("K32").Value = 3155.6316
$ws.Range("K45").Value = 8866.5
$ws.Range("H45").Value = 6957
$ws.Range("I45").Value = 8866.5
$ws.Range("M45").Value = -8489.5
$ws.Range("J134").Value = 99999.5
$ws.Range("L134").Value = 99999.5
$ws.Range("H134").Value = 99999.5
$ws.Range("N134").Value = -110139.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K20").Value = 1154.6364
$ws.Range("M20").Value = -907.6364000000001
$ws.Range("N20").Value = -2259.5714
$ws.Range("I20").Value = 1154.6364
$ws.Range("H20").Value = 1392.2222
$ws.Range("L20").Value = 1765.5714
$ws.Range("J20").Value = 1765.5714
$ws.Range("I26").Value = 28966.334
$ws.Range("H26").Value = 28966.334
$ws.Range("K26").Value = 28966.334
$ws.Range("M26").Value = -28674.334
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("K33").Value = 0
$ws.Range("H33").Value = 7500
$ws.Range("N33").Value = -8172
$ws.Range("M33").ClearContents()
$ws.Range("J33").Value = 7500
$ws.Range("L33").Value = 7500
$ws.Range("I33").Value = 0
$ws.Range("I94").Value = 2466.625
$ws.Range("K94").Value = 2466.625
$ws.Range("H94").Value = 2439.1765
$ws.Range("M94").Value = -2015.625
$ws.Range("J134").Value = 4000
$ws.Range("L134").Value = 12000
$ws.Range("K134").Value = 250006440
$ws.Range("H134").Value = 71430984
$ws.Range("N134").Value = -17070
$ws.Range("I134").Value = 83335480
$ws.Range("M134").Value = -250003905
$ws.Range("M138").Value = -94760
$ws.Range("L138").Value = 99999.5
$ws.Range("J138").Value = 99999.5
$ws.Range("K138").Value = 99900
$ws.Range("I138").Value = 99900
$ws.Range("H138").Value = 99966.336
$ws.Range("N138").Value = -110279.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -6020.615
$ws.Range("J31").Value = 1353.5
$ws.Range("H31").Value = 4425.2856
$ws.Range("L31").Value = 1353.5
$ws.Range("N31").Value = -1943.5
$ws.Range("I31").Value = 6315.615
$ws.Range("K31").Value = 6315.615
$ws.Range("L34").Value = 1353.5
$ws.Range("J34").Value = 1353.5
$ws.Range("H34").Value = 4425.2856
$ws.Range("K34").Value = 6315.615
$ws.Range("N34").Value = -1757.5
$ws.Range("M34").Value = -6113.615
$ws.Range("I34").Value = 6315.615
$ws.Range("M74").ClearContents()
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("H74").Value = 39493
$ws.Range("H77").Value = 39493
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("I77").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K3").Value = 23026.9995
$ws.Range("I3").Value = 7675.6665
$ws.Range("H3").Value = 7675.6665
$ws.Range("M3").Value = -22914.9995
$ws.Range("I18").Value = 847.375
$ws.Range("K18").Value = 2542.125
$ws.Range("H18").Value = 1577.9
$ws.Range("M18").Value = -2373.125
$ws.Range("L62").Value = 0
$ws.Range("H62").Value = 11001
$ws.Range("J62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("J65").Value = 0
$ws.Range("H65").Value = 11001
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("J102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("I102").Value = 0
$ws.Range("H102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("K140").Value = 3072.5295
$ws.Range("M140").Value = 2107.4705
$ws.Range("H140").Value = 1133.9445
$ws.Range("I140").Value = 1024.1765

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K2").Value = 53.75
$ws.Range("M2").Value = 59.25
$ws.Range("H2").Value = 166.35715
$ws.Range("L2").Value = 316.5
$ws.Range("J2").Value = 316.5
$ws.Range("I2").Value = 53.75
$ws.Range("N2").Value = -542.5
$ws.Range("H97").Value = 730.05
$ws.Range("N97").Value = -1947.6
$ws.Range("I97").Value = 654.86664
$ws.Range("M97").Value = -158.86664
$ws.Range("L97").Value = 955.6
$ws.Range("J97").Value = 955.6
$ws.Range("K97").Value = 654.86664
$ws.Range("I122").Value = 4789.625
$ws.Range("N122").Value = -40900
$ws.Range("H122").Value = 6231.7
$ws.Range("J122").Value = 12000
$ws.Range("M122").Value = -11918.875
$ws.Range("K122").Value = 14368.875
$ws.Range("L122").Value = 36000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K7").Value = 3577.3125
$ws.Range("M7").Value = -3465.3125
$ws.Range("H7").Value = 3762.4211
$ws.Range("I7").Value = 3577.3125
$ws.Range("K22").Value = 3937.5
$ws.Range("M22").Value = -3642.5
$ws.Range("I22").Value = 3937.5
$ws.Range("H22").Value = 3350.25
$ws.Range("H27").Value = 3350.25
$ws.Range("I27").Value = 3937.5
$ws.Range("K27").Value = 3937.5
$ws.Range("M27").Value = -3830.5
$ws.Range("H40").Value = 2722.5
$ws.Range("L40").Value = 2850
$ws.Range("J40").Value = 2850
$ws.Range("N40").Value = -3122
$ws.Range("M40").Value = -2544
$ws.Range("I40").Value = 2680
$ws.Range("K40").Value = 2680
$ws.Range("K61").Value = 4063.3572
$ws.Range("I61").Value = 4063.3572
$ws.Range("M61").Value = -3861.3572
$ws.Range("H61").Value = 4025.4
$ws.Range("N68").Value = -3053
$ws.Range("I68").Value = 966
$ws.Range("K68").Value = 966
$ws.Range("J68").Value = 1555
$ws.Range("M68").Value = -217
$ws.Range("H68").Value = 1113.25
$ws.Range("L68").Value = 1555
$ws.Range("I71").Value = 966
$ws.Range("N71").Value = -15263
$ws.Range("J71").Value = 1555
$ws.Range("H71").Value = 1113.25
$ws.Range("K71").Value = 4830
$ws.Range("L71").Value = 7775
$ws.Range("M71").Value = -1086
$ws.Range("J74").Value = 99998.5
$ws.Range("M74").Value = -40230.668
$ws.Range("L74").Value = 99998.5
$ws.Range("I74").Value = 41228.668
$ws.Range("N74").Value = -101994.5
$ws.Range("K74").Value = 41228.668
$ws.Range("H74").Value = 55921.125
$ws.Range("N77").Value = -309979.5
$ws.Range("J77").Value = 99998.5
$ws.Range("L77").Value = 299995.5
$ws.Range("H77").Value = 55921.125
$ws.Range("K77").Value = 123686.004
$ws.Range("M77").Value = -118694.004
$ws.Range("I77").Value = 41228.668
$ws.Range("H82").Value = 1600
$ws.Range("M82").Value = -1139
$ws.Range("I82").Value = 1500
$ws.Range("K82").Value = 1500
$ws.Range("I85").Value = 1500
$ws.Range("K85").Value = 1500
$ws.Range("H85").Value = 1600
$ws.Range("M85").Value = -252
$ws.Range("K113").Value = 4063.3572
$ws.Range("H113").Value = 4025.4
$ws.Range("I113").Value = 4063.3572
$ws.Range("M113").Value = -1893.3572
$ws.Range("I126").Value = 3577.3125
$ws.Range("M126").Value = -8261.9375
$ws.Range("H126").Value = 3762.4211
$ws.Range("K126").Value = 10731.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K2").Value = 50000
$ws.Range("M2").Value = -49888
$ws.Range("H2").Value = 50000
$ws.Range("I2").Value = 50000
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("L8").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J32").Value = 20000
$ws.Range("N32").Value = -20634
$ws.Range("L32").Value = 20000
$ws.Range("H32").Value = 20000
$ws.Range("H34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("I34").Value = 0
$ws.Range("K62").Value = 4485
$ws.Range("I62").Value = 4485
$ws.Range("H62").Value = 5548.5293
$ws.Range("M62").Value = -3861
$ws.Range("M65").Value = -19305
$ws.Range("K65").Value = 22425
$ws.Range("I65").Value = 4485
$ws.Range("H65").Value = 5548.5293
$ws.Range("K113").Value = 2292.1365
$ws.Range("H113").Value = 832.1786
$ws.Range("I113").Value = 764.0454999999999
$ws.Range("M113").Value = -122.1364999999996
$ws.Range("M132").Value = -250004330
$ws.Range("I132").Value = 83335620
$ws.Range("H132").Value = 83335620
$ws.Range("K132").Value = 250006860

Write-Output "Applied 255 cell edits"
